$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Locate & delete the old paragraph text (the two big runs that used
#    to read "If one wants ... brittle" + " and hard.It is true ... a third")
# ---------------------------------------------------------------------
$oldText = "If one wants to have it strong and very fine, one has to gather it before the seed is completely ripe. Because if you wait for the seed to be completely ripe, the base is so dried by the sun that it is brittle and hard.It is true that like this, the seed reduces by a third"

$findRange = $d.Content
$found = $findRange.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    Write-Output "ERROR: could not find old paragraph text"
}

$global:editStart = $findRange.Start
$delRange = $d.Range($findRange.Start, $findRange.End)
$delRange.Text = ""

# The old trailing "." run (previously colored, now to be rebuilt with
# no color as part of the $runs list below) is still sitting right
# after our deletion point - remove it too so we don't end up with a
# duplicate full stop.
$oldDotRange = $d.Range($global:editStart, $global:editStart + 1)
if ($oldDotRange.Text -ne ".") {
    Write-Output "ERROR: expected a leftover '.' run, found instead:"
    Write-Output $oldDotRange.Text
}
$oldDotRange.Text = ""

# $global:pos tracks where the next new run should be inserted.
$global:pos = $global:editStart

# ---------------------------------------------------------------------
# We rebuild the whole paragraph by cloning formatting from existing
# runs elsewhere in the document via Copy/Paste (this reproduces the
# exact minimal <w:rPr> that Word's Font API otherwise refuses to
# reproduce faithfully, e.g. runs with no explicit color, or runs that
# would otherwise pick up stray formatting from an adjacent run). After
# pasting we overwrite the pasted range's .Text with the real content
# we want, which keeps the run's formatting intact.
#
# Anchors used (all exist untouched elsewhere in the document):
#   "dissolved"  -> plain run, <w:color val="000000"/><w:rtl val="0"/>
#   "gum"        -> plain run, <w:rtl val="0"/>                (no color)
#   "amp;"       -> grey Courier run, rFonts+color=a9a9a9+sz18+rtl
#
# We always (re)locate the anchors starting AFTER our edit zone, so we
# never accidentally grab a copy we ourselves just pasted earlier in
# the document (our edits only ever grow the target paragraph, which
# sits before all three anchors in the document).
# ---------------------------------------------------------------------

function Get-SafeSearchStart() {
    $p = $d.Paragraphs.Item(14)
    return $p.Range.End
}

$runs = @(
    @{Text="Who"; Kind="none"},
    @{Text=" wants to have it "; Kind="color"},
    @{Text="beautiful"; Kind="none"},
    @{Text=" "; Kind="color"},
    @{Text="&"; Kind="none"},
    @{Text="amp;"; Kind="courier"},
    @{Text=" very "; Kind="color"},
    @{Text="delicate"; Kind="none"},
    @{Text=", one "; Kind="color"},
    @{Text="needs to"; Kind="none"},
    @{Text=" "; Kind="color"},
    @{Text="pick "; Kind="none"},
    @{Text="it before the seed "; Kind="color"},
    @{Text="may be"; Kind="none"},
    @{Text=" "; Kind="color"},
    @{Text="perfectly"; Kind="none"},
    @{Text=" ripe. "; Kind="color"},
    @{Text="For"; Kind="none"},
    @{Text=" if "; Kind="color"},
    @{Text="one"; Kind="none"},
    @{Text=" waits for the seed to be "; Kind="color"},
    @{Text="perfectly"; Kind="none"},
    @{Text=" ripe, the "; Kind="color"},
    @{Text="foot"; Kind="none"},
    @{Text=" is so dried by the sun that it is brittle"; Kind="color"},
    @{Text=" "; Kind="color"},
    @{Text="&"; Kind="none"},
    @{Text="amp;"; Kind="courier"},
    @{Text=" "; Kind="color"},
    @{Text="breakable"; Kind="none"},
    @{Text=". It is true that "; Kind="color"},
    @{Text="in this way"; Kind="none"},
    @{Text=" the seed is reduce"; Kind="color"},
    @{Text="d"; Kind="none"},
    @{Text=" by a third"; Kind="color"},
    @{Text="."; Kind="none"}
)

foreach ($item in $runs) {
    $txt = $item.Text
    $kind = $item.Kind

    if ($kind -eq "color") {
        $anchorWord = "dissolved"
        $anchorLen = 9
    }
    elseif ($kind -eq "none") {
        $anchorWord = "gum"
        $anchorLen = 3
    }
    else {
        $anchorWord = "amp;"
        $anchorLen = 4
    }

    $searchStart = Get-SafeSearchStart
    $anchor = $d.Range($searchStart, $d.Content.End)
    $ok = $anchor.Find.Execute($anchorWord, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "ERROR: could not find anchor"
        Write-Output $anchorWord
    }
    $anchor.Copy()

    $insPoint = $d.Range($global:pos, $global:pos)
    $insPoint.Paste()

    $pastedRange = $d.Range($global:pos, $global:pos + $anchorLen)
    $pastedRange.Text = $txt

    $global:pos = $global:pos + $txt.Length
}

Write-Output "Final paragraph text:"
Write-Output $d.Paragraphs.Item(14).Range.Text
